$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 696.5
$ws.Range("I28").Value = 390.2
$ws.Range("J28").Value = 1352.8572
$ws.Range("K28").Value = 390.2
$ws.Range("L28").Value = 1352.8572
$ws.Range("M28").Value = 94.80000000000001
$ws.Range("N28").Value = -2322.8572

$ws.Range("H86").Value = 2105.6897
$ws.Range("I86").Value = 2050.9092
$ws.Range("J86").Value = 2277.8572
$ws.Range("K86").Value = 2050.9092
$ws.Range("L86").Value = 2277.8572
$ws.Range("M86").Value = -927.9092000000001
$ws.Range("N86").Value = -4523.8572

$ws.Range("H88").Value = 6812.5
$ws.Range("J88").Value = 8500
$ws.Range("L88").Value = 8500
$ws.Range("N88").Value = -9312

$ws.Range("H89").Value = 2105.6897
$ws.Range("I89").Value = 2050.9092
$ws.Range("J89").Value = 2277.8572
$ws.Range("K89").Value = 10254.546
$ws.Range("L89").Value = 11389.286
$ws.Range("M89").Value = -4638.546
$ws.Range("N89").Value = -22621.286

$ws.Range("H91").Value = 6812.5
$ws.Range("J91").Value = 8500
$ws.Range("L91").Value = 8500
$ws.Range("N91").Value = -11308

$ws.Range("H92").Value = 451.53845
$ws.Range("I92").Value = 414.16666
$ws.Range("K92").Value = 414.16666
$ws.Range("M92").Value = 833.83334

$ws.Range("H98").Value = 610.94116
$ws.Range("I98").Value = 606.3077
$ws.Range("J98").Value = 626
$ws.Range("K98").Value = 606.3077
$ws.Range("L98").Value = 626
$ws.Range("M98").Value = 891.6923
$ws.Range("N98").Value = -3622

$ws.Range("H118").Value = 10958
$ws.Range("I118").Value = 15011.429
$ws.Range("K118").Value = 45034.287
$ws.Range("M118").Value = -43377.287

$ws.Range("H122").Value = 610.94116
$ws.Range("I122").Value = 606.3077
$ws.Range("J122").Value = 626
$ws.Range("K122").Value = 1818.9231
$ws.Range("L122").Value = 1878
$ws.Range("M122").Value = 631.0769
$ws.Range("N122").Value = -6778

$ws.Range("H128").Value = 26000
$ws.Range("I128").Value = 26000
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 26000
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -21020
$ws.Range("N128").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 28060.648
$ws.Range("I2").Value = 1145.25
$ws.Range("J2").Value = 111797.445
$ws.Range("K2").Value = 1145.25
$ws.Range("L2").Value = 111797.445
$ws.Range("M2").Value = -1032.25
$ws.Range("N2").Value = -112023.445

$ws.Range("H32").Value = 24595.082
$ws.Range("I32").Value = 5862.671
$ws.Range("J32").Value = 166961.4
$ws.Range("K32").Value = 5862.671
$ws.Range("L32").Value = 166961.4
$ws.Range("M32").Value = -5575.671
$ws.Range("N32").Value = -167535.4

$ws.Range("H74").Value = 721.68085
$ws.Range("I74").Value = 680.4872
$ws.Range("J74").Value = 922.5
$ws.Range("K74").Value = 680.4872
$ws.Range("L74").Value = 922.5
$ws.Range("M74").Value = 193.5128
$ws.Range("N74").Value = -2670.5

$ws.Range("H77").Value = 721.68085
$ws.Range("I77").Value = 680.4872
$ws.Range("J77").Value = 922.5
$ws.Range("K77").Value = 3402.436
$ws.Range("L77").Value = 4612.5
$ws.Range("M77").Value = 965.5639999999999
$ws.Range("N77").Value = -13348.5

$ws.Range("H116").Value = 28060.648
$ws.Range("I116").Value = 1145.25
$ws.Range("J116").Value = 111797.445
$ws.Range("K116").Value = 1145.25
$ws.Range("L116").Value = 111797.445
$ws.Range("M116").Value = 1148.75
$ws.Range("N116").Value = -116385.445

$ws.Range("H122").Value = 1615.2759
$ws.Range("I122").Value = 1542.1
$ws.Range("J122").Value = 1777.8889
$ws.Range("K122").Value = 4626.299999999999
$ws.Range("L122").Value = 5333.6667
$ws.Range("M122").Value = -2176.299999999999
$ws.Range("N122").Value = -10233.6667

$ws.Range("H132").Value = 11909.437
$ws.Range("I132").Value = 14076.86
$ws.Range("J132").Value = 4142.8335
$ws.Range("K132").Value = 42230.58
$ws.Range("L132").Value = 12428.5005
$ws.Range("M132").Value = -39700.58
$ws.Range("N132").Value = -17488.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 28060.648
$ws.Range("I3").Value = 1145.25
$ws.Range("J3").Value = 111797.445
$ws.Range("K3").Value = 1145.25
$ws.Range("L3").Value = 111797.445
$ws.Range("M3").Value = -1031.25
$ws.Range("N3").Value = -112025.445

$ws.Range("H86").Value = 86647.69500000001
$ws.Range("I86").Value = 124226.664
$ws.Range("J86").Value = 2095
$ws.Range("K86").Value = 124226.664
$ws.Range("L86").Value = 2095
$ws.Range("M86").Value = -123103.664
$ws.Range("N86").Value = -4341

$ws.Range("H89").Value = 86647.69500000001
$ws.Range("I89").Value = 124226.664
$ws.Range("J89").Value = 2095
$ws.Range("K89").Value = 621133.3200000001
$ws.Range("L89").Value = 10475
$ws.Range("M89").Value = -615517.3200000001
$ws.Range("N89").Value = -21707

$ws.Range("H99").Value = 2253.76
$ws.Range("I99").Value = 1824.1428
$ws.Range("J99").Value = 2420.8333
$ws.Range("K99").Value = 1824.1428
$ws.Range("L99").Value = 2420.8333
$ws.Range("M99").Value = -326.1428000000001
$ws.Range("N99").Value = -5416.8333

$ws.Range("H107").Value = 76923920
$ws.Range("I107").Value = 100000820
$ws.Range("J107").Value = 933.3333
$ws.Range("K107").Value = 100000820
$ws.Range("L107").Value = 933.3333
$ws.Range("M107").Value = -99998900
$ws.Range("N107").Value = -4773.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2233.8
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 2233.8
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 6701.400000000001
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -8323.400000000001

$ws.Range("H71").Value = 2233.8
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 2233.8
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 20104.2
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -28216.2

$ws.Range("H131").Value = 2127.6626
$ws.Range("J131").Value = 2263.0134
$ws.Range("L131").Value = 6789.040199999999
$ws.Range("N131").Value = -16869.0402

$ws.Range("H140").Value = 6915.5557
$ws.Range("I140").Value = 7480
$ws.Range("J140").Value = 2400
$ws.Range("K140").Value = 22440
$ws.Range("L140").Value = 7200
$ws.Range("M140").Value = -17260
$ws.Range("N140").Value = -17560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 3500
$ws.Range("J18").Value = 3500
$ws.Range("L18").Value = 3500
$ws.Range("N18").Value = -4086

$ws.Range("H43").Value = 5229.75
$ws.Range("I43").Value = 1900
$ws.Range("J43").Value = 6339.6665
$ws.Range("K43").Value = 1900
$ws.Range("L43").Value = 6339.6665
$ws.Range("M43").Value = -1749
$ws.Range("N43").Value = -6641.6665

$ws.Range("H70").Value = 45730.6
$ws.Range("I70").Value = 60480.332
$ws.Range("J70").Value = 7802.7144
$ws.Range("K70").Value = 60480.332
$ws.Range("L70").Value = 7802.7144
$ws.Range("M70").Value = -60210.332
$ws.Range("N70").Value = -8342.714400000001

$ws.Range("H73").Value = 45730.6
$ws.Range("I73").Value = 60480.332
$ws.Range("J73").Value = 7802.7144
$ws.Range("K73").Value = 60480.332
$ws.Range("L73").Value = 7802.7144
$ws.Range("M73").Value = -59544.332
$ws.Range("N73").Value = -9674.714400000001

$ws.Range("H122").Value = 1595.3939
$ws.Range("I122").Value = 1373.7142
$ws.Range("J122").Value = 1983.3334
$ws.Range("K122").Value = 4121.142599999999
$ws.Range("L122").Value = 5950.0002
$ws.Range("M122").Value = -1671.142599999999
$ws.Range("N122").Value = -10850.0002

$ws.Range("H132").Value = 2402.318
$ws.Range("I132").Value = 1865.921
$ws.Range("J132").Value = 5799.5
$ws.Range("K132").Value = 5597.763
$ws.Range("L132").Value = 17398.5
$ws.Range("M132").Value = -3067.763
$ws.Range("N132").Value = -22458.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 19500
$ws.Range("J74").Value = 19500
$ws.Range("L74").Value = 19500
$ws.Range("N74").Value = -21496

$ws.Range("H77").Value = 19500
$ws.Range("J77").Value = 19500
$ws.Range("L77").Value = 58500
$ws.Range("N77").Value = -68484

$ws.Range("H93").Value = 1896.7667
$ws.Range("I93").Value = 1799.2106
$ws.Range("J93").Value = 2065.2727
$ws.Range("K93").Value = 1799.2106
$ws.Range("L93").Value = 2065.2727
$ws.Range("M93").Value = -551.2106000000001
$ws.Range("N93").Value = -4561.2727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 167379.83
$ws.Range("I81").Value = 700
$ws.Range("J81").Value = 250719.75
$ws.Range("K81").Value = 1400
$ws.Range("L81").Value = 501439.5
$ws.Range("M81").Value = -339
$ws.Range("N81").Value = -503561.5

$ws.Range("H84").Value = 167379.83
$ws.Range("I84").Value = 700
$ws.Range("J84").Value = 250719.75
$ws.Range("K84").Value = 7000
$ws.Range("L84").Value = 2507197.5
$ws.Range("M84").Value = -1696
$ws.Range("N84").Value = -2517805.5

$ws.Range("H107").Value = 2000
$ws.Range("I107").Value = 2000
$ws.Range("K107").Value = 6000
$ws.Range("M107").Value = -4080
